$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subscript-3 unicode character used in the PEPE price (0.0₃0512)
$sub3 = [char]0x2083

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.669.89"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +4.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.050.54"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +3.61%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.44"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.70"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +6.67%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.048.84"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +3.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.505"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +5.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.156"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +8.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.10"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -4.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.479"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +10.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000231"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +6.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.88"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.548.94"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.762.34"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +4.32%  "
$ws.Range("E17").Value = "  +3.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.051.19"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +3.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.73"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "479.17"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.11"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +7.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.677"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +6.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.58"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +9.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.09"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +15.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.26"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +4.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.79"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +4.76%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.96"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +7.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.04"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +5.25%  "
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.23"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +4.88%  "
$ws.Range("E32").Value = "  +2.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.44"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +6.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.66"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +4.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.20"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +8.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.88"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0408"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +7.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "439.36"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0810"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +3.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.84"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +22.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.966.69"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.23"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +5.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.113"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "27.96"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +7.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.260"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +7.57%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.13"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +10.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.113"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +5.55%  "
$ws.Range("B49").Value = "PEPE"
$ws.Range("C49").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$d49 = "0.0{0}0512" -f $sub3
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = $d49
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +7.16%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "116.89"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.08"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +6.90%  "
